# NpcSkill.xlsx - "Skill" sheet data refresh
# 5 monster skill sets completed: No.107, No.108, Jake_B, Jake_R, Doncina(01/02)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# Wipe the previous (now-stale) table before writing the new layout.
$ws.Range("A1:G10").ClearContents()

# --- Header row ---
$ws.Cells.Item(1,1).Value = "cid"
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "nameKor"
$ws.Cells.Item(1,4).Value = "styleTube"
$ws.Cells.Item(1,5).Value = "enhancerTube"
$ws.Cells.Item(1,6).Value = "coolerTube"
$ws.Cells.Item(1,7).Value = "relicTube"

# --- Row 2: No.107 ---
$ws.Cells.Item(2,1).Value = 19100
$ws.Cells.Item(2,2).Value = "19_deadlyattack"
$ws.Cells.Item(2,3).Value = "데들리 어택"
$ws.Cells.Item(2,4).Value = "no107_style"
$ws.Cells.Item(2,5).Value = "no107_enhancer"
$ws.Cells.Item(2,6).Value = "no107_cooler"

# --- Row 3: No.108 ---
$ws.Cells.Item(3,1).Value = 19101
$ws.Cells.Item(3,2).Value = "19_criticalhit"
$ws.Cells.Item(3,3).Value = "크리티컬 히트"
$ws.Cells.Item(3,4).Value = "no108_style"
$ws.Cells.Item(3,5).Value = "no108_enhancer"
$ws.Cells.Item(3,6).Value = "no108_cooler"

# --- Row 4: Jake_B ---
$ws.Cells.Item(4,1).Value = 19103
$ws.Cells.Item(4,2).Value = "19_jakebounce"
$ws.Cells.Item(4,3).Value = "썩은폭탄"
$ws.Cells.Item(4,4).Value = "jake_B_style"
$ws.Cells.Item(4,5).Value = "jake_B_enhancer"
$ws.Cells.Item(4,6).Value = "jake_B_cooler"

# --- Row 5: Jake_R ---
$ws.Cells.Item(5,1).Value = 19104
$ws.Cells.Item(5,2).Value = "19_jakerange"
$ws.Cells.Item(5,3).Value = "썩은돌팔매"
$ws.Cells.Item(5,4).Value = "jake_R_style"
$ws.Cells.Item(5,5).Value = "jake_R_enhancer"
$ws.Cells.Item(5,6).Value = "jake_R_cooler"

# --- Row 6: Doncina01 ---
$ws.Cells.Item(6,1).Value = 19200
$ws.Cells.Item(6,2).Value = "19_MustleMustle"
$ws.Cells.Item(6,3).Value = "머슬머슬"
$ws.Cells.Item(6,4).Value = "doncina01_style"
$ws.Cells.Item(6,5).Value = "doncina01_enhancer"
$ws.Cells.Item(6,6).Value = "doncina01_cooler"

# --- Row 7: Doncina02 ---
$ws.Cells.Item(7,1).Value = 19201
$ws.Cells.Item(7,2).Value = "19_HustleHustle"
$ws.Cells.Item(7,3).Value = "허슬허슬"
$ws.Cells.Item(7,4).Value = "doncina02_style"
$ws.Cells.Item(7,5).Value = "doncina02_enhancer"
$ws.Cells.Item(7,6).Value = "doncina02_cooler"
$ws.Cells.Item(7,7).Value = "doncina02_relic"

# --- Column widths (bestFit-style, matches the authored layout) ---
$ws.Columns.Item(4).ColumnWidth = 14.71
$ws.Columns.Item(5).ColumnWidth = 19.29
$ws.Columns.Item(6).ColumnWidth = 16.39
$ws.Columns.Item(7).ColumnWidth = 14.41

# --- View state: zoomed to 85%, selection parked on C10 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("C10").Select()
